$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Add the two new calculated columns to the table (Tabela1) ---
$colG = $lo.ListColumns.Add()
$ws.Range("G1").Value = "Implementação"

$colH = $lo.ListColumns.Add()
$ws.Range("H1").Value = "Teste"

# --- Halve the "TAMANHO GLOBAL" multiplier in column E (was *3, now *1.5) ---
# Row 2 holds its own (non-shared) formula.
$ws.Range("E2").Formula = "=D2*1.5"

# Rows 3:36 share one formula group (master in E3); re-assigning the whole
# range keeps them grouped as a shared formula, same as the original file.
$ws.Range("E3:E36").Formula = "=D3*1.5"

# --- Populate the new Implementação / Teste columns for every data row ---
# Each row gets its own individual formula (not a shared-formula group),
# matching how the source workbook stores these two new columns.
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=Tabela1[[#This Row],[TAMANHO GLOBAL]]*0.5"
    $ws.Cells.Item($r, 8).Formula = "=Tabela1[[#This Row],[TAMANHO GLOBAL]]*0.1"
}

# --- Update the view: selection moves to H2, scroll back to the top row ---
$null = $ws.Range("H2").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
